$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5210496154209789
$ws.Range("D2").Value = 0.6075376934226

$ws.Range("C3").Value = 0.360446983781771
$ws.Range("D3").Value = 0.7219505599865363

$ws.Range("C4").Value = -0.2461234121835621
$ws.Range("D4").Value = 0.8078684987726339

$ws.Range("C5").Value = -0.4149893607479724
$ws.Range("D5").Value = 0.6821688173759157

$ws.Range("C6").Value = -0.2142552586919674
$ws.Range("D6").Value = 0.8323214302350395

$ws.Range("C7").Value = -0.7005275699367309
$ws.Range("D7").Value = 0.4909409182285236

$ws.Range("C8").Value = -0.9912173772625391
$ws.Range("D8").Value = 0.3323569975550367

$ws.Range("C9").Value = -0.6184229030711621
$ws.Range("D9").Value = 0.5426469887349041

$ws.Range("C10").Value = -0.6008909314992762
$ws.Range("D10").Value = 0.5540537040917608

$ws.Range("C11").Value = -0.3120915228638611
$ws.Range("D11").Value = 0.7579096767446589
